$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03034166666666667
$ws.Range("H2").Value = 0.09102499999999999
$ws.Range("I2").Value = 0.001862306623420098
$ws.Range("J2").Value = 0.001862306623420098
$ws.Range("M2").Value = 8.521337333333333
$ws.Range("N2").Value = 25.564012
$ws.Range("O2").Value = 0.2943426187002489
$ws.Range("P2").Value = 0.2943426187002489
$ws.Range("Q2").Value = 0.2585515769222222
$ws.Range("R2").Value = 2.3269641923
$ws.Range("S2").Value = 0.0005481562083602898
$ws.Range("T2").Value = 0.0005481562083602898

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03034166666666667
$ws.Range("H3").Value = 0.09102499999999999
$ws.Range("I3").Value = 0.001862306623420098
$ws.Range("J3").Value = 0.001862306623420098
$ws.Range("O3").Value = 0.1683364841626613
$ws.Range("P3").Value = 0.1683364841626613
$ws.Range("Q3").Value = 0.1478673514083333
$ws.Range("R3").Value = 1.330806162675
$ws.Range("S3").Value = 0.0003134941494193765
$ws.Range("T3").Value = 0.0003134941494193765

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03034166666666667
$ws.Range("H4").Value = 0.09102499999999999
$ws.Range("I4").Value = 0.001862306623420098
$ws.Range("J4").Value = 0.001862306623420098
$ws.Range("O4").Value = 0.5373208971370899
$ws.Range("P4").Value = 0.53732089713709
$ws.Range("Q4").Value = 0.4719845392472222
$ws.Range("R4").Value = 4.247860853224999
$ws.Range("S4").Value = 0.001000656265640431
$ws.Range("T4").Value = 0.001000656265640432

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.5257850852532363
$ws.Range("J5").Value = 0.5257850852532362
$ws.Range("M5").Value = 8.521337333333333
$ws.Range("N5").Value = 25.564012
$ws.Range("O5").Value = 0.2943426187002489
$ws.Range("P5").Value = 0.2943426187002489
$ws.Range("Q5").Value = 72.99687452367689
$ws.Range("R5").Value = 656.971870713092
$ws.Range("S5").Value = 0.1547609588669712
$ws.Range("T5").Value = 0.1547609588669711

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.5257850852532363
$ws.Range("J6").Value = 0.5257850852532362
$ws.Range("O6").Value = 0.1683364841626613
$ws.Range("P6").Value = 0.1683364841626613
$ws.Range("S6").Value = 0.08850881267669493
$ws.Range("T6").Value = 0.08850881267669491

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.5257850852532363
$ws.Range("J7").Value = 0.5257850852532362
$ws.Range("O7").Value = 0.5373208971370899
$ws.Range("P7").Value = 0.53732089713709
$ws.Range("S7").Value = 0.2825153137095702
$ws.Range("T7").Value = 0.2825153137095702

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 7.695814000000001
$ws.Range("I8").Value = 0.4723526081233437
$ws.Range("J8").Value = 0.4723526081233436
$ws.Range("M8").Value = 8.521337333333333
$ws.Range("N8").Value = 25.564012
$ws.Range("O8").Value = 0.2943426187002489
$ws.Range("P8").Value = 0.2943426187002489
$ws.Range("Q8").Value = 65.57862714858935
$ws.Range("R8").Value = 590.207644337304
$ws.Range("S8").Value = 0.1390335036249174
$ws.Range("T8").Value = 0.1390335036249174

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 7.695814000000001
$ws.Range("I9").Value = 0.4723526081233437
$ws.Range("J9").Value = 0.4723526081233436
$ws.Range("O9").Value = 0.1683364841626613
$ws.Range("P9").Value = 0.1683364841626613
$ws.Range("Q9").Value = 37.504849209926
$ws.Range("S9").Value = 0.07951417733654699
$ws.Range("T9").Value = 0.079514177336547

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 7.695814000000001
$ws.Range("I10").Value = 0.4723526081233437
$ws.Range("J10").Value = 0.4723526081233436
$ws.Range("O10").Value = 0.5373208971370899
$ws.Range("P10").Value = 0.53732089713709
$ws.Range("S10").Value = 0.2538049271618792
$ws.Range("T10").Value = 0.2538049271618793
